$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting (bold, border, centered) from the last existing label row (A1376) down to the new label rows
$ws.Range("A1376").Copy() | Out-Null
$ws.Range("A1377:A1488").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$labels = @(
  'pop division emp_okev at 0 in 2020',
  'pop division emp_okev at 0 in 2025',
  'pop division emp_okev at א-טור in 2020',
  'pop division emp_okev at א-טור in 2025',
  'pop division emp_okev at אבו טור סילוא in 2020',
  'pop division emp_okev at אבו טור סילוא in 2025',
  'pop division emp_okev at אזור תעשיה גב in 2020',
  'pop division emp_okev at אזור תעשיה גב in 2025',
  'pop division emp_okev at אזור תעשיה עט in 2020',
  'pop division emp_okev at אזור תעשיה עט in 2025',
  'pop division emp_okev at בית הכרם in 2020',
  'pop division emp_okev at בית הכרם in 2025',
  'pop division emp_okev at בית צפפא,שרפת in 2020',
  'pop division emp_okev at בית צפפא,שרפת in 2025',
  'pop division emp_okev at בקעה מושבות in 2020',
  'pop division emp_okev at בקעה מושבות in 2025',
  'pop division emp_okev at גבעת המטוס in 2020',
  'pop division emp_okev at גבעת המטוס in 2025',
  'pop division emp_okev at גבעת משואה in 2020',
  'pop division emp_okev at גבעת משואה in 2025',
  'pop division emp_okev at גילה in 2020',
  'pop division emp_okev at גילה in 2025',
  'pop division emp_okev at הגבעה הצרפתי in 2020',
  'pop division emp_okev at הגבעה הצרפתי in 2025',
  'pop division emp_okev at העיר העתיקה in 2020',
  'pop division emp_okev at העיר העתיקה in 2025',
  'pop division emp_okev at הקריות של ירו in 2020',
  'pop division emp_okev at הקריות של ירו in 2025',
  'pop division emp_okev at הר החוצבים in 2020',
  'pop division emp_okev at הר החוצבים in 2025',
  'pop division emp_okev at הר הצופים in 2020',
  'pop division emp_okev at הר הצופים in 2025',
  'pop division emp_okev at הר נוף in 2020',
  'pop division emp_okev at הר נוף in 2025',
  'pop division emp_okev at חומת שמואל - ה in 2020',
  'pop division emp_okev at חומת שמואל - ה in 2025',
  'pop division emp_okev at חומת שמואל-מע in 2020',
  'pop division emp_okev at חומת שמואל-מע in 2025',
  'pop division emp_okev at טלביה רחביה in 2020',
  'pop division emp_okev at טלביה רחביה in 2025',
  'pop division emp_okev at יובלים גנים in 2020',
  'pop division emp_okev at יובלים גנים in 2025',
  'pop division emp_okev at כניסה לעיר in 2020',
  'pop division emp_okev at כניסה לעיר in 2025',
  'pop division emp_okev at כניסה לעיר - ק in 2020',
  'pop division emp_okev at כניסה לעיר - ק in 2025',
  'pop division emp_okev at כפר עקב in 2020',
  'pop division emp_okev at כפר עקב in 2025',
  'pop division emp_okev at מורשה in 2020',
  'pop division emp_okev at מורשה in 2025',
  'pop division emp_okev at מחנה שועפט in 2020',
  'pop division emp_okev at מחנה שועפט in 2025',
  'pop division emp_okev at מנחת - מלחה in 2020',
  'pop division emp_okev at מנחת - מלחה in 2025',
  'pop division emp_okev at מעלות דפנה in 2020',
  'pop division emp_okev at מעלות דפנה in 2025',
  'pop division emp_okev at מער in 2020',
  'pop division emp_okev at מער in 2025',
  'pop division emp_okev at מער חרדי in 2020',
  'pop division emp_okev at מער חרדי in 2025',
  'pop division emp_okev at מער מזרח ואדי in 2020',
  'pop division emp_okev at מער מזרח ואדי in 2025',
  'pop division emp_okev at מקור ברוך in 2020',
  'pop division emp_okev at מקור ברוך in 2025',
  'pop division emp_okev at מרכז רפואי הד in 2020',
  'pop division emp_okev at מרכז רפואי הד in 2025',
  'pop division emp_okev at נוה יעקב in 2020',
  'pop division emp_okev at נוה יעקב in 2025',
  'pop division emp_okev at נחלאות שערי ח in 2020',
  'pop division emp_okev at נחלאות שערי ח in 2025',
  'pop division emp_okev at ניות -גבעת מר in 2020',
  'pop division emp_okev at ניות -גבעת מר in 2025',
  'pop division emp_okev at סנהדריה רמת א in 2020',
  'pop division emp_okev at סנהדריה רמת א in 2025',
  'pop division emp_okev at עין כרם in 2020',
  'pop division emp_okev at עין כרם in 2025',
  'pop division emp_okev at עיסאוויה in 2020',
  'pop division emp_okev at עיסאוויה in 2025',
  'pop division emp_okev at פסגת זאב in 2020',
  'pop division emp_okev at פסגת זאב in 2025',
  'pop division emp_okev at פת - קטמונים in 2020',
  'pop division emp_okev at פת - קטמונים in 2025',
  'pop division emp_okev at צור באהר מורח in 2020',
  'pop division emp_okev at צור באהר מורח in 2025',
  'pop division emp_okev at קטמון הישנה in 2020',
  'pop division emp_okev at קטמון הישנה in 2025',
  'pop division emp_okev at קריית הממשלה- in 2020',
  'pop division emp_okev at קריית הממשלה- in 2025',
  'pop division emp_okev at קריית משה -גב in 2020',
  'pop division emp_okev at קריית משה -גב in 2025',
  'pop division emp_okev at ראס אל עמד - ג'' in 2020',
  'pop division emp_okev at ראס אל עמד - ג'' in 2025',
  'pop division emp_okev at רכס לבן in 2020',
  'pop division emp_okev at רכס לבן in 2025',
  'pop division emp_okev at רמות in 2020',
  'pop division emp_okev at רמות in 2025',
  'pop division emp_okev at רמת שלמה in 2020',
  'pop division emp_okev at רמת שלמה in 2025',
  'pop division emp_okev at רמת שרת - בית in 2020',
  'pop division emp_okev at רמת שרת - בית in 2025',
  'pop division emp_okev at רסקו in 2020',
  'pop division emp_okev at רסקו in 2025',
  'pop division emp_okev at שדה התעופה עט in 2020',
  'pop division emp_okev at שדה התעופה עט in 2025',
  'pop division emp_okev at שועפאט - בית in 2020',
  'pop division emp_okev at שועפאט - בית in 2025',
  'pop division emp_okev at תלפיות in 2020',
  'pop division emp_okev at תלפיות in 2025',
  'pop division emp_okev at תלפיות ארנונ in 2020',
  'pop division emp_okev at תלפיות ארנונ in 2025',
  'pop division emp_okev at תלפיות מזרח in 2020',
  'pop division emp_okev at תלפיות מזרח in 2025'
)

$cvals = @(
  10.61797177433101,
  48.66136791876972,
  8.879335743536435,
  8.97809595681324,
  10.48817622266559,
  10.64197920075806,
  0,
  0,
  0,
  0,
  7.124781198522624,
  7.308352879444294,
  11.71001852449868,
  11.90347984974643,
  22.00711043041655,
  22.15309767411908,
  0,
  0,
  5.585780140371284,
  6.280991735537189,
  5.351811824539101,
  5.728047129564301,
  39.23725427102401,
  40.49165572956781,
  0,
  0,
  20.70118699601915,
  24.91431893195359,
  0,
  0,
  0,
  0,
  14.66547261650109,
  15.1941214789956,
  6.944442962067281,
  6.942148760330578,
  0,
  0,
  7.347476667204568,
  6.738142506191168,
  5.202732404265584,
  5.692485246231844,
  0,
  0,
  91.24879455885682,
  90.56489963902756,
  12.08791208791209,
  11.78516286509618,
  0,
  0,
  11.74716806238538,
  11.44824005909651,
  6.832876037424109,
  7.767185784048269,
  18.77653961963095,
  19.44024615580287,
  0,
  0,
  32.70664157449798,
  35.5682969497193,
  0,
  0,
  0,
  0,
  0,
  0,
  16.39725505146778,
  18.45514042445038,
  102.0295600503881,
  113.8681090059341,
  5.945316652078143,
  6.942148760330579,
  13.74233004967821,
  16.76729304575607,
  4.505540946656641,
  4.753474830954171,
  10.99870325494772,
  11.16113328814083,
  7.457743220880475,
  7.623295056673723,
  5.811428131537058,
  6.172487281127843,
  9.504076162580008,
  9.552561591347882,
  3.956323603879162,
  4.171014397564975,
  0,
  0,
  221.3329223895845,
  237.4078039230274,
  10.47478667238044,
  10.38830448361933,
  0,
  0,
  19.00032971832463,
  19.4693353446468,
  16.63789607831781,
  19.27171326363419,
  10.1945169955856,
  11.15606071240911,
  4.423537535334355,
  4.318931062344308,
  25.71537651982733,
  28.19310963582926,
  9.893383537176318,
  9.956893577718436,
  0,
  0,
  9.000384485519868,
  9.360846884237608,
  4.816265316393286,
  5.21737476075624
)

$dvals = @(
  10.61797177433101,
  50.26030791605782,
  8.879335743536435,
  8.97809595681324,
  10.48817622266559,
  10.65733287027282,
  0,
  0,
  0,
  0,
  7.124781198522624,
  6.891837625228979,
  11.71001852449868,
  11.90347984974643,
  22.00711043041655,
  22.95406507450567,
  0,
  0,
  5.585780140371284,
  6.28099173553719,
  5.351811824539101,
  5.72773368218265,
  39.23725427102401,
  40.94104047252055,
  0,
  0,
  20.70118699601915,
  24.91280982354486,
  0,
  0,
  0,
  0,
  14.66547261650109,
  15.1941214789956,
  6.944442962067281,
  6.942148760330582,
  0,
  0,
  7.347476667204568,
  6.717886164423779,
  5.202732404265584,
  5.689438851389011,
  0,
  0,
  91.24879455885682,
  91.28950915822267,
  12.08791208791209,
  11.78516286509618,
  0,
  0,
  11.74716806238538,
  11.44704172880582,
  6.832876037424109,
  7.767185784048263,
  18.77653961963095,
  19.31312751089166,
  0,
  0,
  32.70664157449798,
  35.34646507018998,
  0,
  0,
  0,
  0,
  0,
  0,
  16.39725505146778,
  18.49033157702685,
  102.0295600503881,
  115.1956358320139,
  5.945316652078143,
  6.942148760330581,
  13.74233004967821,
  16.670043633783,
  4.505540946656641,
  4.75347483095417,
  10.99870325494772,
  11.16113328814083,
  7.457743220880475,
  7.64009813006245,
  5.811428131537058,
  6.099233051633656,
  9.504076162580008,
  9.552347929157344,
  3.956323603879162,
  4.172972218938142,
  0,
  0,
  221.3329223895845,
  230.1613518100752,
  10.47478667238044,
  10.38663976374064,
  0,
  0,
  19.00032971832463,
  19.37580411701953,
  16.63789607831781,
  19.28558597677667,
  10.1945169955856,
  11.04595970731102,
  4.423537535334355,
  4.312667561891613,
  25.71537651982733,
  27.03759665161701,
  9.893383537176318,
  9.956915770717734,
  0,
  0,
  9.000384485519868,
  9.382152700650103,
  4.816265316393286,
  5.207898905781394
)

$evals = @(
  10.61797177433101,
  49.7684289776571,
  8.879335743536435,
  8.978095956813242,
  10.48817622266559,
  10.65123788087779,
  0,
  0,
  0,
  0,
  7.124781198522624,
  6.961116780234643,
  11.71001852449868,
  11.90347984974643,
  22.00711043041655,
  22.69873846913899,
  0,
  0,
  5.585780140371284,
  6.28099173553719,
  5.351811824539101,
  5.727903274549295,
  39.23725427102401,
  40.93406273939863,
  0,
  0,
  20.70118699601915,
  24.91152165368718,
  0,
  0,
  0,
  0,
  14.66547261650109,
  15.1941214789956,
  6.944442962067281,
  6.942148760330578,
  0,
  0,
  7.347476667204568,
  6.718213713655519,
  5.202732404265584,
  5.68618828405399,
  0,
  0,
  91.24879455885682,
  91.39560998488896,
  12.08791208791209,
  11.78516286509618,
  0,
  0,
  11.74716806238538,
  11.44705978651874,
  6.832876037424109,
  7.767185784048269,
  18.77653961963095,
  19.32225336826691,
  0,
  0,
  32.70664157449798,
  35.39803257083106,
  0,
  0,
  0,
  0,
  0,
  0,
  16.39725505146778,
  18.48548096986664,
  102.0295600503881,
  118.8240836730507,
  5.945316652078143,
  6.942148760330581,
  13.74233004967821,
  16.49860679659169,
  4.505540946656641,
  4.753474830954172,
  10.99870325494772,
  11.16113328814084,
  7.457743220880475,
  7.636071415732046,
  5.811428131537058,
  6.119577814703307,
  9.504076162580008,
  9.552350672556532,
  3.956323603879162,
  4.172428755729251,
  0,
  0,
  221.3329223895845,
  234.523447818702,
  10.47478667238044,
  10.38646121365346,
  0,
  0,
  19.00032971832463,
  19.24772766594151,
  16.63789607831781,
  19.28381182603322,
  10.1945169955856,
  10.92233902608439,
  4.423537535334355,
  4.313317859992234,
  25.71537651982733,
  26.71540738503945,
  9.893383537176318,
  9.956923310275558,
  0,
  0,
  9.000384485519868,
  9.374972223665662,
  4.816265316393286,
  5.208887551814326
)

$startRow = 1377
for ($i = 0; $i -lt $labels.Length; $i++) {
  $r = $startRow + $i
  $ws.Cells.Item($r, 1).Value = $labels[$i]
  $ws.Cells.Item($r, 3).Value = $cvals[$i]
  $ws.Cells.Item($r, 4).Value = $dvals[$i]
  $ws.Cells.Item($r, 5).Value = $evals[$i]
}

Write-Output "Added $($labels.Length) rows"